# Apply the commit's data update: the three species-observation records
# that were stored in rows 3, 4 and 6 get rotated:
#   new row 3 <- old row 4 data
#   new row 4 <- old row 6 data
#   new row 6 <- old row 3 data
# (Rows 2, 5, 7 and all other columns are untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3 (receives what used to be row 4's record) ----
$ws.Range("A3").Value = 111835718
$ws.Range("B3").Value = 56398
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
$ws.Range("J3").ClearContents()
$ws.Range("M3").Value = "äldre spår"
$ws.Range("Q3").Value = 471101.0270993827
$ws.Range("R3").Value = 6810411.753755242
$ws.Range("S3").Value = 10
$ws.Range("AF3").ClearContents()

# ---- Row 4 (receives what used to be row 6's record) ----
$ws.Range("A4").Value = 111835826
$ws.Range("Q4").Value = 470915.776864712
$ws.Range("R4").Value = 6810385.536630718
$ws.Range("S4").Value = 5
$ws.Range("AC4").Value = "även hackspettbo, troligen av tret"

# ---- Row 6 (receives what used to be row 3's record) ----
$ws.Range("A6").Value = 111835758
$ws.Range("B6").Value = 77550
$ws.Range("E6").Value = 185
$ws.Range("F6").Value = "Violettgrå tagellav"
$ws.Range("G6").Value = "Bryoria nadvornikiana"
$ws.Range("H6").Value = "(Gyeln.) Brodo & D.Hawksw."
$ws.Range("L6").ClearContents()
$ws.Range("M6").ClearContents()
$ws.Range("Q6").Value = 471087.4311846643
$ws.Range("R6").Value = 6810390.807424263
$ws.Range("S6").Value = 5
$ws.Range("AC6").ClearContents()
